{"js": "// Remove the \"File/Directory (Non Java)\" section: its Heading1 title\n// paragraph plus the two body paragraphs describing file/folder naming\n// rules (the \"Coding rules.docx\" / \"documenti\" examples).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text.trim();\n  if (\n    text === \"File/Directory (Non Java)\" ||\n    text.indexOf(\"I nomi dei file devono iniziare con una lettera maiuscola\") === 0 ||\n    text.indexOf(\"I nomi delle cartelle devono iniziare con una lettera minuscola\") === 0\n  ) {\n    targets.push(p);\n  }\n}\n\nfor (const p of targets) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"File/Directory (Non Java)\" section: its Heading1 title\n# paragraph plus the two body paragraphs describing the file/folder naming\n# rules (the \"Coding rules.docx\" / \"documenti\" examples).\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$startIdx = -1\n$endIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($startIdx -lt 0 -and $t -eq \"File/Directory (Non Java)\") {\n        $startIdx = $i\n        $endIdx = $i\n    }\n    elseif ($startIdx -ge 0 -and $endIdx -eq $startIdx) {\n        if ($t.StartsWith(\"I nomi dei file devono iniziare\")) {\n            $endIdx = $i\n        }\n    }\n    elseif ($startIdx -ge 0 -and $endIdx -gt $startIdx) {\n        if ($t.StartsWith(\"I nomi delle cartelle devono iniziare\")) {\n            $endIdx = $i\n        }\n    }\n}\n\nif ($startIdx -gt 0 -and $endIdx -ge $startIdx) {\n    $startPara = $d.Paragraphs.Item($startIdx)\n    $endPara = $d.Paragraphs.Item($endIdx)\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
